# Katalog guncellendi - Pzt 24.11.2025 15:40:47,11
# Adds 4 new "YELEK" (vest) products to the product catalog sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New product rows: urun_adi, fiyat, kategori, gorsel, aciklama, stok
$yeniUrunler = @(
    @("ŞİŞME YELEK 5020 SOFT HAKİ", "400 TL", "YELEK", "5020HAKİ.jpg"),
    @("ŞİŞME YELEK 5020 SOFT  KAHVERENGİ", "400 TL", "YELEK", "5020KAHVERENGİ.jpg"),
    @("ŞİŞME YELEK 5020 SOFT LACİVERT", "400 TL", "YELEK", "5020LACİVERT.jpg"),
    @("ŞİŞME YELEK 5020 SOFT SİYAH", "400 TL", "YELEK", "5020SİYAH.jpg")
)

$aciklama = "S-M-L-XL-2XL Beden seçeneği mevcuttur.Ürünümüz serili olarak satılmaktadır.Belirtilen fiyatlar adet fiyatıdır."
$stok = "Var"

$satir = 76
foreach ($urun in $yeniUrunler) {
    $ws.Cells.Item($satir, 1).Value = $urun[0]
    $ws.Cells.Item($satir, 2).Value = $urun[1]
    $ws.Cells.Item($satir, 3).Value = $urun[2]
    $ws.Cells.Item($satir, 4).Value = $urun[3]
    $ws.Cells.Item($satir, 5).Value = $aciklama
    $ws.Cells.Item($satir, 6).Value = $stok
    $satir = $satir + 1
}

$ws.Range("B79").Select()
$excel.ActiveWindow.ScrollRow = 55
$excel.ActiveWindow.ScrollColumn = 1
